{"js": "const body = context.document.body;\n\n// --- 1) \"Versi\" + \"on\" -> single run \"Version\" ---------------------------\n// search() finds \"Version\" across the two existing runs. insertText(...,\n// \"Replace\") on that hit rewrites/merges the underlying runs. Because the\n// replacement text is identical to the current text (\"Version\"), route\n// through a distinct placeholder first so the rewrite actually happens,\n// then rename the placeholder back to \"Version\".\nlet results = body.search(\"Version\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst placeholder = \"VersionPLACEHOLDER\";\nresults.items[0].insertText(placeholder, \"Replace\");\nawait context.sync();\n\nresults = body.search(placeholder, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n\n// --- 2) \" 2\" -> \" 1.\" -----------------------------------------------------\nresults = body.search(\"2\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"1.\", \"Replace\");\nawait context.sync();\n\n// --- 3) drop the now-redundant trailing \".\" run ---------------------------\n// (originally a separate run placed after the _GoBack bookmark). After step\n// 2 there are two \".\" matches in the paragraph (\"...1.\" and the old trailing\n// \".\"); the last occurrence is the original stray run.\nresults = body.search(\".\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[results.items.length - 1].delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) \"Versi\" + \"on\" -> single run \"Version\" -------------------------------\n# The two runs are merged into one whenever Word rewrites a Range spanning\n# both of them with text that differs from the current text. Since the\n# target text (\"Version\") is identical to the current text, first stamp a\n# distinct placeholder (forces the merge/rewrite) and then rename it back.\n$verRange = $d.Content\n$verRange.Find.ClearFormatting()\n$verRange.Find.Text = \"Version\"\n$null = $verRange.Find.Execute()\n$verStart = $verRange.Start\n\n$placeholder = \"VersionPLACEHOLDER\"\n$verRange.Text = $placeholder\n$mergedRange = $d.Range($verStart, $verStart + $placeholder.Length)\n$mergedRange.Text = \"Version\"\n\n# --- 2) \" 2\" -> \" 1.\" ----------------------------------------------------\n$numRange = $d.Content\n$numRange.Find.ClearFormatting()\n$numRange.Find.Text = \"2\"\n$null = $numRange.Find.Execute()\n$numRange.Text = \"1.\"\n\n# --- 3) drop the now-redundant trailing \".\" run --------------------------\n# (originally a separate run placed after the _GoBack bookmark)\n$paraRange = $d.Paragraphs(1).Range\n$paraEnd = $paraRange.End\n$lastChar = $d.Range($paraEnd - 2, $paraEnd - 1)\nif ($lastChar.Text -eq \".\") {\n    $lastChar.Delete()\n}\n"}
